# Workbook/sheet handles (already open per the runtime contract).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row: a 4th population label added below the existing three
# (Pt0_blast / Pt0_Er1 / Pt0_NK) -> adds "Pt0_Bcells" as a new shared
# string and writes it into A4.
$ws.Range("A4").Value = "Pt0_Bcells"

# Re-apply formatting to the two rows above it (cosmetic "touch" seen in
# the saved file as a style-index bump on A2/A3).
$ws.Range("A2:A3").Style = "Normal"

# Selection moved from A7 to C3.
$ws.Range("C3").Select()
